# Applies the commit "Updated symbol list on Thu Jan 19 09:23:40 UTC 2023 with GitHub Actions"
# Refreshes Price (D) and Volume(1h) (E) columns on Sheet1 for the crypto symbol list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'291.24"
$ws.Range("E2").Value = "'-3.22%"

$ws.Range("D3").Value = "'30.68"
$ws.Range("E3").Value = "'-6.37%"

$ws.Range("D4").Value = "'4.955"
$ws.Range("E4").Value = "'0.29%"

$ws.Range("D5").Value = "'0.07211"
$ws.Range("E5").Value = "'-6.41%"

$ws.Range("D6").Value = "'1.818"
$ws.Range("E6").Value = "'-7.79%"

$ws.Range("D7").Value = "'7.691"
$ws.Range("E7").Value = "'-1.85%"

$ws.Range("E8").Value = "'-0.79%"

$ws.Range("D9").Value = "'0.8960"
$ws.Range("E9").Value = "'-2.59%"

$ws.Range("D10").Value = "'0.1656"
$ws.Range("E10").Value = "'-5.34%"

$ws.Range("D11").Value = "'0.07716"
$ws.Range("E11").Value = "'-0.59%"

$ws.Range("D12").Value = "'0.08097"
$ws.Range("E12").Value = "'-5.89%"

$ws.Range("D13").Value = "'0.03022"
$ws.Range("E13").Value = "'-4.43%"

$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'-0.33%"

$ws.Range("D15").Value = "'0.001502"
$ws.Range("E15").Value = "'-0.57%"

$ws.Range("D16").Value = "'0.005699"
$ws.Range("E16").Value = "'-3.80%"

$ws.Range("D18").Value = "'3.466"
$ws.Range("E18").Value = "'0.22%"

$ws.Range("D19").Value = "'2.080"
$ws.Range("E19").Value = "'-3.39%"

$ws.Range("E20").Value = "'-0.92%"

$ws.Range("E21").Value = "'-2.08%"

$ws.Range("D22").Value = "'4.038"
$ws.Range("E22").Value = "'-5.83%"

$ws.Range("D23").Value = "'0.2389"
$ws.Range("E23").Value = "'19.93%"

$ws.Range("D24").Value = "'0.04501"
$ws.Range("E24").Value = "'-0.35%"

$ws.Range("E25").Value = "'-0.86%"

$ws.Range("D26").Value = "'0.004009"
$ws.Range("E26").Value = "'-9.16%"

$ws.Range("D27").Value = "'0.0001251"
$ws.Range("E27").Value = "'0.03%"

$ws.Range("D39").Value = "'0.01592"
$ws.Range("E39").Value = "'-6.73%"

$ws.Range("D40").Value = "'0.04405"
$ws.Range("E40").Value = "'-5.99%"

$ws.Range("D41").Value = "'0.007290"
$ws.Range("E41").Value = "'-2.47%"

$ws.Range("D42").Value = "'0.009953"

$ws.Range("D43").Value = "'0.1306"
$ws.Range("E43").Value = "'-3.44%"

$ws.Range("E44").Value = "'-11.99%"

$ws.Range("D45").Value = "'0.009517"
$ws.Range("E45").Value = "'-9.23%"

$ws.Range("D46").Value = "'0.00005953"
$ws.Range("E46").Value = "'-4.91%"

$ws.Range("E47").Value = "'0.05%"

$ws.Range("E48").Value = "'173.65%"

$ws.Range("D49").Value = "'0.003003"
$ws.Range("E49").Value = "'-3.29%"

$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.05%"

$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.05%"

Write-Host "Updated 68 cells across the price/volume columns."
